# Updates the cryptocurrency price/volume table (and, for rows 39-40,
# the coin name + link) to the refreshed scrape values.
# Each value is written with a leading apostrophe so Excel stores it as
# text (matching the original inline-string cells) rather than silently
# re-typing numeric-looking strings (e.g. "219.56") as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.427.08"
$ws.Range("E2").Value = "'  +0.35%  "

$ws.Range("D3").Value = "'1.701.48"
$ws.Range("E3").Value = "'  +0.95%  "

$ws.Range("E4").Value = "'  +0.06%  "

$ws.Range("D5").Value = "'219.56"
$ws.Range("E5").Value = "'  +0.62%  "

$ws.Range("D6").Value = "'0.5513"
$ws.Range("E6").Value = "'  +5.37%  "

$ws.Range("E7").Value = "'  +0.00%  "

$ws.Range("D8").Value = "'0.2746"
$ws.Range("E8").Value = "'  +1.77%  "

$ws.Range("D9").Value = "'0.06479"
$ws.Range("E9").Value = "'  +1.18%  "

$ws.Range("D10").Value = "'22.12"
$ws.Range("E10").Value = "'  +0.58%  "

$ws.Range("D11").Value = "'0.07683"
$ws.Range("E11").Value = "'  +2.56%  "

$ws.Range("D12").Value = "'1.698.32"
$ws.Range("E12").Value = "'  -1.05%  "

$ws.Range("D13").Value = "'4.557"
$ws.Range("E13").Value = "'  -0.07%  "

$ws.Range("D14").Value = "'0.5851"
$ws.Range("E14").Value = "'  +1.12%  "

$ws.Range("D15").Value = "'0.000008400"
$ws.Range("E15").Value = "'  -0.84%  "

$ws.Range("D16").Value = "'65.83"
$ws.Range("E16").Value = "'  +2.46%  "

$ws.Range("D17").Value = "'26.455.74"

$ws.Range("D18").Value = "'4.955"
$ws.Range("E18").Value = "'  +0.74%  "

$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "'  +0.13%  "

$ws.Range("E20").Value = "'  +1.14%  "

$ws.Range("D21").Value = "'192.32"
$ws.Range("E21").Value = "'  +2.10%  "

$ws.Range("D22").Value = "'6.265"
$ws.Range("E22").Value = "'  +1.24%  "

$ws.Range("E23").Value = "'  +0.00%  "

$ws.Range("D24").Value = "'148.96"
$ws.Range("E24").Value = "'  +3.13%  "

$ws.Range("D25").Value = "'0.1329"
$ws.Range("E25").Value = "'  +8.09%  "

$ws.Range("D26").Value = "'7.943"
$ws.Range("E26").Value = "'  +3.08%  "

$ws.Range("D27").Value = "'15.86"
$ws.Range("E27").Value = "'  +0.50%  "

$ws.Range("D28").Value = "'0.06309"
$ws.Range("E28").Value = "'  -4.79%  "

$ws.Range("D29").Value = "'1.383"
$ws.Range("E29").Value = "'  +2.45%  "

$ws.Range("D30").Value = "'1.333"
$ws.Range("E30").Value = "'  +0.30%  "

$ws.Range("D31").Value = "'3.614"
$ws.Range("E31").Value = "'  +1.24%  "

$ws.Range("D32").Value = "'3.616"
$ws.Range("E32").Value = "'  +1.40%  "

$ws.Range("E33").Value = "'  +1.92%  "

$ws.Range("E34").Value = "'  +1.99%  "

$ws.Range("D35").Value = "'0.6191"
$ws.Range("E35").Value = "'  -0.22%  "

$ws.Range("E36").Value = "'  +0.31%  "

$ws.Range("D37").Value = "'2.742"
$ws.Range("E37").Value = "'  +1.69%  "

$ws.Range("D38").Value = "'0.01651"
$ws.Range("E38").Value = "'  +2.24%  "

$ws.Range("B39").Value = "'Maker"
$ws.Range("C39").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "'1.120.94"
$ws.Range("E39").Value = "'  +0.91%  "

$ws.Range("B40").Value = "'FraxShare"
$ws.Range("C40").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.184"
$ws.Range("E40").Value = "'  -2.92%  "

$ws.Range("D41").Value = "'0.8830"
$ws.Range("E41").Value = "'  +0.27%  "

$ws.Range("E42").Value = "'  +0.00%  "

$ws.Range("D43").Value = "'101.33"
$ws.Range("E43").Value = "'  +0.23%  "

$ws.Range("D44").Value = "'1.850.27"
$ws.Range("E44").Value = "'  +0.81%  "

$ws.Range("D45").Value = "'57.83"
$ws.Range("E45").Value = "'  +1.90%  "

$ws.Range("D46").Value = "'0.00000000107"
$ws.Range("E46").Value = "'  -3.28%  "

$ws.Range("D47").Value = "'8.242"
$ws.Range("E47").Value = "'  +0.94%  "

$ws.Range("E48").Value = "'  +0.07%  "

$ws.Range("D49").Value = "'0.05275"

$ws.Range("D50").Value = "'6.126"
$ws.Range("E50").Value = "'  +1.16%  "

$ws.Range("D51").Value = "'0.4305"
$ws.Range("E51").Value = "'  -0.05%  "
